# Update the "想去人数" (wanted-to-go count) column (F) on the 展览 and 全部类型
# sheets to reflect the latest scrape numbers.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 1726
    4  = 793
    6  = 36
    7  = 12024
    10 = 480
    11 = 415
    12 = 1115
    13 = 870
    14 = 13489
    15 = 13506
    17 = 155
    20 = 490
    21 = 97
    23 = 1899
    24 = 177
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
